$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Remove the date stamp that was previously placed in C1 (value + date
# format), restoring the sheet to only use columns A:B (dimension A1:B24).
$ws.Range("C1").Clear() | Out-Null
